$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '68.711.50'
$ws.Cells.Item(2, 5).Value = '  +2.17%  '

$ws.Cells.Item(3, 4).Value = '3.921.69'
$ws.Cells.Item(3, 5).Value = '  +1.13%  '

$ws.Cells.Item(4, 5).Value = '  +0.06%  '

$ws.Cells.Item(5, 4).Value = "'481.29"
$ws.Cells.Item(5, 5).Value = '  +1.99%  '

$ws.Cells.Item(6, 4).Value = "'145.19"
$ws.Cells.Item(6, 5).Value = '  +0.00%  '

$ws.Cells.Item(7, 4).Value = "'0.622"
$ws.Cells.Item(7, 5).Value = '  -1.80%  '

$ws.Cells.Item(8, 4).Value = "'0.997"
$ws.Cells.Item(8, 5).Value = '  -0.19%  '

$ws.Cells.Item(9, 5).Value = '  -2.85%  '

$ws.Cells.Item(10, 5).Value = '  +7.90%  '

$ws.Cells.Item(11, 4).Value = "'0.0000353"
$ws.Cells.Item(11, 5).Value = '  +11.66%  '

$ws.Cells.Item(12, 4).Value = "'42.69"
$ws.Cells.Item(12, 5).Value = '  -1.90%  '

$ws.Cells.Item(13, 2).Value = 'Polkadot'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(13, 4).Value = "'10.51"
$ws.Cells.Item(13, 5).Value = '  +0.86%  '

$ws.Cells.Item(14, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(14, 4).Value = '4.555.06'
$ws.Cells.Item(14, 5).Value = '  +0.81%  '

$ws.Cells.Item(15, 2).Value = 'Uniswap'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(15, 4).Value = "'14.60"
$ws.Cells.Item(15, 5).Value = '  -1.47%  '

$ws.Cells.Item(16, 2).Value = 'WrappedEther'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(16, 4).Value = '3.916.99'
$ws.Cells.Item(16, 5).Value = '  +0.78%  '

$ws.Cells.Item(17, 5).Value = '  -0.31%  '

$ws.Cells.Item(18, 5).Value = '  -2.03%  '

$ws.Cells.Item(19, 5).Value = '  -2.58%  '

$ws.Cells.Item(20, 4).Value = '68.737.10'
$ws.Cells.Item(20, 5).Value = '  +1.78%  '

$ws.Cells.Item(21, 4).Value = "'433.23"
$ws.Cells.Item(21, 5).Value = '  -0.41%  '

$ws.Cells.Item(22, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(22, 4).Value = "'14.61"
$ws.Cells.Item(22, 5).Value = '  -1.62%  '

$ws.Cells.Item(23, 2).Value = 'ImmutableX'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(23, 4).Value = "'3.36"
$ws.Cells.Item(23, 5).Value = '  +2.69%  '

$ws.Cells.Item(24, 4).Value = "'87.53"
$ws.Cells.Item(24, 5).Value = '  -1.69%  '

$ws.Cells.Item(25, 4).Value = "'11.67"
$ws.Cells.Item(25, 5).Value = '  +17.23%  '

$ws.Cells.Item(26, 5).Value = '  -1.44%  '

$ws.Cells.Item(27, 4).Value = "'38.17"
$ws.Cells.Item(27, 5).Value = '  +0.58%  '

$ws.Cells.Item(28, 5).Value = '  +7.78%  '

$ws.Cells.Item(29, 4).Value = "'10.17"
$ws.Cells.Item(29, 5).Value = '  -1.06%  '

$ws.Cells.Item(30, 4).Value = "'711.77"
$ws.Cells.Item(30, 5).Value = '  -2.44%  '

$ws.Cells.Item(31, 2).Value = 'Hedera'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(31, 4).Value = "'0.130"
$ws.Cells.Item(31, 5).Value = '  -3.66%  '

$ws.Cells.Item(32, 2).Value = 'Cosmos'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(32, 4).Value = "'13.27"
$ws.Cells.Item(32, 5).Value = '  -3.98%  '

$ws.Cells.Item(33, 4).Value = "'2.86"
$ws.Cells.Item(33, 5).Value = '  +3.50%  '

$ws.Cells.Item(34, 4).Value = '0.0₃0916'
$ws.Cells.Item(34, 5).Value = '  +33.00%  '

$ws.Cells.Item(35, 4).Value = "'41.38"
$ws.Cells.Item(35, 5).Value = '  -5.40%  '

$ws.Cells.Item(36, 4).Value = "'58.53"
$ws.Cells.Item(36, 5).Value = '  +1.09%  '

$ws.Cells.Item(37, 5).Value = '  -7.76%  '

$ws.Cells.Item(38, 5).Value = '  +2.62%  '

$ws.Cells.Item(39, 5).Value = '  -0.02%  '

$ws.Cells.Item(40, 5).Value = '  -2.34%  '

$ws.Cells.Item(41, 4).Value = "'3.05"
$ws.Cells.Item(41, 5).Value = '  +9.29%  '

$ws.Cells.Item(42, 4).Value = "'2.76"
$ws.Cells.Item(42, 5).Value = '  +6.48%  '

$ws.Cells.Item(43, 5).Value = '  +2.84%  '

$ws.Cells.Item(44, 5).Value = '  -2.94%  '

$ws.Cells.Item(45, 5).Value = '  -0.21%  '

$ws.Cells.Item(46, 5).Value = '  -0.05%  '

$ws.Cells.Item(47, 5).Value = '  -1.32%  '

$ws.Cells.Item(48, 5).Value = '  -0.72%  '

$ws.Cells.Item(49, 4).Value = "'147.41"
$ws.Cells.Item(49, 5).Value = '  +2.05%  '

$ws.Cells.Item(50, 5).Value = '  -4.36%  '

$ws.Cells.Item(51, 5).Value = '  -1.88%  '
